$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @("MACD", "RSI", "Signal_line", "VIX_short", "close_short", "close_long", "VIX", "VIX_long", "fedrate")
$importances = @(0.3930817313426748, 0.3080310985219448, 0.1320233135593225, 0.04050396027363317, 0.03403665931781309, 0.03071072798458908, 0.02369486540704165, 0.02211115837664275, 0.01580648521633805)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $features[$i]
    $ws.Cells.Item($row, 2).Value = $importances[$i]
}
